# Add new columns I (I0) and J (IF) to the weekly data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers "I0" and "IF" ---
$ws.Cells.Item(1, 9).Value  = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Match the header formatting used by the existing header cells (B1:H1):
# bold font, centered horizontally, top-aligned vertically, thin border all round.
foreach ($colIdx in @(9, 10)) {
    $hdr = $ws.Cells.Item(1, $colIdx)
    $hdr.Font.Bold = $true
    $hdr.HorizontalAlignment = -4108   # xlCenter
    $hdr.VerticalAlignment = -4160     # xlTop
    $hdr.Borders.LineStyle = 1         # xlContinuous (thin)
}

# --- Data rows 2-33: values for columns I (I0) and J (IF) ---
$data = @(
    @(2, 10, 10),
    @(3, 14, 18),
    @(4, 1, 6),
    @(5, 1, 6),
    @(6, 1, 5),
    @(7, 1, 5),
    @(8, 1, 7),
    @(9, 1, 5),
    @(10, 1, 5),
    @(11, 1, 7),
    @(12, 1, 6),
    @(13, 1, 7),
    @(14, 1, 6),
    @(15, 1, 5),
    @(16, 1, 6),
    @(17, 1, 7),
    @(18, 1, 6),
    @(19, 1, 6),
    @(20, 1, 6),
    @(21, 1, 5),
    @(22, 1, 7),
    @(23, 1, 7),
    @(24, 1, 5),
    @(25, 1, 4),
    @(26, 1, 5),
    @(27, 1, 7),
    @(28, 1, 6),
    @(29, 1, 6),
    @(30, 1, 6),
    @(31, 1, 4),
    @(32, 1, 4),
    @(33, 1, 2)
)

foreach ($row in $data) {
    $r = $row[0]
    $i0 = $row[1]
    $iF = $row[2]
    $ws.Cells.Item($r, 9).Value  = $i0
    $ws.Cells.Item($r, 10).Value = $iF
}

Write-Output "I0 and IF columns added"
